$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.745571851730347
$ws.Range("B1").Value = 2.629175901412964
$ws.Range("C1").Value = 3.221799373626709
$ws.Range("D1").Value = 1.226133465766907
$ws.Range("E1").Value = 0.8155626058578491
